# Actualización 11 de Mayo - Mañana
# Rebuild the "Rescatables" table with the refreshed roster (20 students,
# rows 2:21) and matching grupo/totales data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$data = @(
    @(18330051920253, "CONTRERAS", "CASTRO", "SAMUEL", "TEMAS DE FILOSOFÍA", "6ALCM", 2),
    @(18330051920281, "LOPEZ", "ROSAS", "ALEXA", "TEMAS DE FILOSOFÍA", "6ALCM", 2),
    @(18330051920309, "VALENCIA", "TORRES", "LUIS ENRIQUE", "TEMAS DE FILOSOFÍA", "6ALCM", 2),
    @(18330051920367, "GASCA", "RUIZ", "ESTEFANI JHOSSELYNE", "TEMAS DE FILOSOFÍA", "6APM", 2),
    @(18330051920368, "GIL", "BANDALA", "AELEN", "TEMAS DE FILOSOFÍA", "6ALCM", 2),
    @(18330051920217, "MARQUEZ", "HERNANDEZ", "ANDRES", "TEMAS DE FILOSOFÍA", "6ARHM", 2),
    @(18330051920329, "GUERRA", "ROMERO", "JOCELYN", "TEMAS DE FILOSOFÍA", "6BLCM", 2),
    @(18330051920424, "LOPEZ", "APALE", "MARIA LARET", "TEMAS DE FILOSOFÍA", "6ALCM", 2),
    @(18330051920339, "MORALES", "TREJO", "ROCIO TAMARA", "TEMAS DE FILOSOFÍA", "6ALCM", 2),
    @(18330051920342, "REYES", "SARMIENTO", "INGRID PAOLA", "TEMAS DE FILOSOFÍA", "6ALCM", 2),
    @(18330051920347, "TEXOCO", "DE JESUS", "MAYTE", "TEMAS DE FILOSOFÍA", "6ALCM", 2),
    @(18330051920352, "ZEPAHUA", "JUAREZ", "QUETZALI", "TEMAS DE FILOSOFÍA", "6ALCM", 2),
    @(18330051920261, "DE GANTE", "GUTIERREZ", "DANIELA", "TEMAS DE FILOSOFÍA", "6ALCM", 1),
    @(18330051920302, "SANCHEZ", "REYES", "YAIRA GUADALUPE", "TEMAS DE FILOSOFÍA", "6ALCM", 1),
    @(18330051920308, "VALDES", "MARIN", "BRANDON", "TEMAS DE FILOSOFÍA", "6ALCM", 1),
    @(18330051920311, "VASQUEZ", "PIEDRAS", "MARTHA MARISOL", "TEMAS DE FILOSOFÍA", "6ALCM", 1),
    @(18330051920440, "DURAN", "CORTES", "LUIS GABRIEL", "TEMAS DE FILOSOFÍA", "6APM", 1),
    @(18330051920364, "ESPINOZA", "AREVALO", "MARY JOSE", "TEMAS DE FILOSOFÍA", "6APM", 1),
    @(18330051920366, "GARCIA", "ARENAS", "JOSUA", "TEMAS DE FILOSOFÍA", "6APM", 1),
    @(18330051920333, "LOPEZ", "MORALES", "SILVIA ESMERALDA", "TEMAS DE FILOSOFÍA", "6BLCM", 1)
)

$rowCount = $data.Length
$colCount = 7
$startRow = 2

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $colCount; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}

Write-Output "Rescatables updated: $rowCount rows written"
